# error solve ifrs list
# Update the numeric financial figures for data rows 2-6 (company rows 1-5)
# and clear out the figures in rows 7-9 (company rows 6-8), leaving only the
# identifying columns A, B and C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---------------------------------------------------------------
$ws.Range("D2").Value = 2198
$ws.Range("E2").Value = 89
$ws.Range("F2").Value = 89
$ws.Range("G2").Value = 84
$ws.Range("H2").Value = 68
$ws.Range("I2").Value = 68
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1838
$ws.Range("L2").Value = 712
$ws.Range("M2").Value = 1126
$ws.Range("N2").Value = 1112
$ws.Range("O2").Value = 15
$ws.Range("P2").Value = 80
$ws.Range("Q2").Value = 28
$ws.Range("R2").Value = -123
$ws.Range("S2").Value = 56
$ws.Range("T2").Value = 73
$ws.Range("U2").Value = -44
$ws.Range("V2").Value = 429
$ws.Range("W2").Value = 4.06
$ws.Range("X2").Value = 3.1
$ws.Range("Y2").Value = 6.29
$ws.Range("Z2").Value = 3.87
$ws.Range("AA2").Value = 63.25
$ws.Range("AB2").Value = 1274.68
$ws.Range("AC2").Value = 423
$ws.Range("AD2").Value = 8.65
$ws.Range("AE2").Value = 6948
$ws.Range("AF2").Value = 0.53
$ws.Range("AG2").Value = 30
$ws.Range("AH2").Value = 0.82
$ws.Range("AI2").Value = 7.08
$ws.Range("AJ2").Value = 16000000

# --- Row 3 ---------------------------------------------------------------
$ws.Range("D3").Value = 2148
$ws.Range("E3").Value = 127
$ws.Range("F3").Value = 127
$ws.Range("G3").Value = 127
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 99
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1773
$ws.Range("L3").Value = 550
$ws.Range("M3").Value = 1223
$ws.Range("N3").Value = 1207
$ws.Range("O3").Value = 15
$ws.Range("P3").Value = 80
$ws.Range("Q3").Value = 282
$ws.Range("R3").Value = -134
$ws.Range("S3").Value = -115
$ws.Range("T3").Value = 68
$ws.Range("U3").Value = 215
$ws.Range("V3").Value = 321
$ws.Range("W3").Value = 5.93
$ws.Range("X3").Value = 4.65
$ws.Range("Y3").Value = 8.57
$ws.Range("Z3").Value = 5.53
$ws.Range("AA3").Value = 45.01
$ws.Range("AB3").Value = 1391.4
$ws.Range("AC3").Value = 621
$ws.Range("AD3").Value = 6.58
$ws.Range("AE3").Value = 7546
$ws.Range("AF3").Value = 0.54
$ws.Range("AG3").Value = 50
$ws.Range("AH3").Value = 1.22
$ws.Range("AI3").Value = 8.05
$ws.Range("AJ3").Value = 16000000

# --- Row 4 ---------------------------------------------------------------
$ws.Range("D4").Value = 2011
$ws.Range("E4").Value = 134
$ws.Range("F4").Value = 134
$ws.Range("G4").Value = 147
$ws.Range("H4").Value = 107
$ws.Range("I4").Value = 107
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1957
$ws.Range("L4").Value = 622
$ws.Range("M4").Value = 1335
$ws.Range("N4").Value = 1320
$ws.Range("O4").Value = 15
$ws.Range("P4").Value = 80
$ws.Range("Q4").Value = 105
$ws.Range("R4").Value = -125
$ws.Range("S4").Value = 48
$ws.Range("T4").Value = 64
$ws.Range("U4").Value = 41
$ws.Range("V4").Value = 359
$ws.Range("W4").Value = 6.64
$ws.Range("X4").Value = 5.34
$ws.Range("Y4").Value = 8.5
$ws.Range("Z4").Value = 5.76
$ws.Range("AA4").Value = 46.6
$ws.Range("AB4").Value = 1515.2
$ws.Range("AC4").Value = 672
$ws.Range("AD4").Value = 13.85
$ws.Range("AE4").Value = 8252
$ws.Range("AF4").Value = 1.13
$ws.Range("AG4").Value = 50
$ws.Range("AH4").Value = 0.54
$ws.Range("AI4").Value = 7.44
$ws.Range("AJ4").Value = 16000000

# --- Row 5 ---------------------------------------------------------------
$ws.Range("D5").Value = 2264
$ws.Range("E5").Value = 139
$ws.Range("F5").Value = 139
$ws.Range("G5").Value = 136
$ws.Range("H5").Value = 110
$ws.Range("I5").Value = 111
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2283
$ws.Range("L5").Value = 874
$ws.Range("M5").Value = 1409
$ws.Range("N5").Value = 1395
$ws.Range("O5").Value = 13
$ws.Range("P5").Value = 80
$ws.Range("Q5").Value = 64
$ws.Range("R5").Value = -280
$ws.Range("S5").Value = 183
$ws.Range("T5").Value = 316
$ws.Range("U5").Value = -251
$ws.Range("V5").Value = 533
$ws.Range("W5").Value = 6.14
$ws.Range("X5").Value = 4.87
$ws.Range("Y5").Value = 8.15
$ws.Range("Z5").Value = 5.2
$ws.Range("AA5").Value = 62.04
$ws.Range("AB5").Value = 1640.6
$ws.Range("AC5").Value = 692
$ws.Range("AD5").Value = 7.47
$ws.Range("AE5").Value = 8721
$ws.Range("AF5").Value = 0.59
$ws.Range("AG5").Value = 50
$ws.Range("AH5").Value = 0.97
$ws.Range("AI5").Value = 7.22
$ws.Range("AJ5").Value = 16000000

# --- Row 6 (no J6/O6 in source) ------------------------------------------
$ws.Range("D6").Value = 2360
$ws.Range("E6").Value = 99
$ws.Range("F6").Value = 99
$ws.Range("G6").Value = 97
$ws.Range("H6").Value = 80
$ws.Range("I6").Value = 79
$ws.Range("K6").Value = 2250
$ws.Range("L6").Value = 776
$ws.Range("M6").Value = 1475
$ws.Range("N6").Value = 1460
$ws.Range("P6").Value = 80
$ws.Range("Q6").Value = 46
$ws.Range("R6").Value = -30
$ws.Range("S6").Value = -36
$ws.Range("T6").Value = 112
$ws.Range("U6").Value = -66
$ws.Range("V6").Value = 506
$ws.Range("W6").Value = 4.2
$ws.Range("X6").Value = 3.4
$ws.Range("Y6").Value = 5.54
$ws.Range("Z6").Value = 3.54
$ws.Range("AA6").Value = 52.59
$ws.Range("AB6").Value = 1727
$ws.Range("AC6").Value = 494
$ws.Range("AD6").Value = 8.07
$ws.Range("AE6").Value = 9127
$ws.Range("AF6").Value = 0.44
$ws.Range("AG6").Value = 50
$ws.Range("AH6").Value = 1.25
$ws.Range("AI6").Value = 10.12
$ws.Range("AJ6").Value = 16000000

# --- Rows 7-9: clear all figures, keep only A/B/C identifying columns ----
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
